# Apply the two edits described by the diff:
#  1. Split the ATLAS_URI connection-string run into two runs (same
#     formatting) -- the first keeps "ATLAS_URI = mongodb+srv:" and a new
#     run carries the rest of the (now placeholder-ized) connection string.
#  2. Remove the trailing paragraph that holds only a manual page break,
#     right before the final section properties.

$d = $word.ActiveDocument

# --- Edit 1: split the ATLAS_URI run -----------------------------------
$rng = $d.Content
$old = "ATLAS_URI = mongodb+srv://dbuser:dbuser@cluster0.zthpw.mongodb.net/retokendb?retryWrites=true&w=majority"
$rng.Find.Execute($old)

$splitStart = $rng.Start
$splitPoint = $splitStart + 24   # length of "ATLAS_URI = mongodb+srv:"
$tailEnd = $rng.End

$newTail = "//mongodb+srv://<username>:<password>@cluster0.zthpw.mongodb.net/<mydatabasename>?retryWrites=true&w=majority"
$tailRange = $d.Range($splitPoint, $tailEnd)
$tailRange.Text = $newTail

# Re-apply the (unchanged) run formatting to the new tail text so it lands
# in its own run instead of being silently merged back into the first one.
$tailRange2 = $d.Range($splitPoint, $splitPoint + $newTail.Length)
$tailRange2.Font.Color = 255
$tailRange2.Font.Color = 2697256   # 0x282829 in BGR order, matches w:color 282829

# --- Edit 2: drop the trailing page-break-only paragraph ---------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Delete()
